# Add new "Code Smell" rows (19-22) to the tracking sheet, matching the
# author's commit "Add code smells to xlsx".
#
# The write order below is chosen so the workbook's shared-string table
# gets new unique strings appended in the same order they appear in the
# target file (Magic Number, BoardServiceImplTests, TaskServiceImplTests,
# "Global variables are introduced to replace the magic numbers",
# TaskController, "Introduce variables to decompose the long statements.").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A + B first for the three "Magic Number" rows (19-21)
$ws.Range("A19").Value2 = "Magic Number"
$ws.Range("B19").Value2 = "BoardServiceImplTests"

$ws.Range("A20").Value2 = "Magic Number"
$ws.Range("B20").Value2 = "TaskServiceImplTests"

$ws.Range("A21").Value2 = "Magic Number"
$ws.Range("B21").Value2 = "UserServiceImplTests"

# Column C ("Refactored?") for those same three rows
$ws.Range("C19").Value2 = "Yes"
$ws.Range("C20").Value2 = "Yes"
$ws.Range("C21").Value2 = "Yes"

# Column D (reasoning) for those same three rows - identical text
$ws.Range("D19").Value2 = "Global variables are introduced to replace the magic numbers"
$ws.Range("D20").Value2 = "Global variables are introduced to replace the magic numbers"
$ws.Range("D21").Value2 = "Global variables are introduced to replace the magic numbers"

# Match the style (explicit black font) already used on rows 13-18
$ws.Range("A19:D21").Font.Color = 0

# Final row (22) - "Long Statement" / TaskController
$ws.Range("A22").Value2 = "Long Statement"
$ws.Range("B22").Value2 = "TaskController"
$ws.Range("C22").Value2 = "Yes"
$ws.Range("D22").Value2 = "Introduce variables to decompose the long statements."
$ws.Range("A22:D22").Font.Color = 0

# Match the author's final selection state (top-left/active cell)
$ws.Range("D10").Select()
